# Fill in the missing "Genotype" (column F) values for the 5CKO / RA
# sample group (rows 82-91), which was left blank while every other
# group already had its genotype recorded.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

for ($r = 82; $r -le 91; $r++) {
    $ws.Cells.Item($r, 6).Value = "WT"
}

# Leave the selection where the last edit landed, matching the
# author's on-screen state after making the change.
$excel.ActiveWindow.ScrollRow = 81
$ws.Range("F91").Select()
